$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(204).Insert()

$ws.Cells.Item(204, 1).Value = 10
$ws.Cells.Item(204, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(204, 3).Value = 'La Araucanía'
$ws.Cells.Item(204, 4).Value = 45075
$ws.Cells.Item(204, 5).Value = 9
$ws.Cells.Item(204, 6).Value = 100112031
$ws.Cells.Item(204, 7).Value = 'Poroto verde'
$ws.Cells.Item(204, 8).Value = 'Sin especificar'
$ws.Cells.Item(204, 9).Value = 'Primera'
$ws.Cells.Item(204, 10).Value = 125
$ws.Cells.Item(204, 11).Value = 25000
$ws.Cells.Item(204, 12).Value = 25000
$ws.Cells.Item(204, 13).Value = 25000
$ws.Cells.Item(204, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(204, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(204, 16).Value = 1000
$ws.Cells.Item(204, 17).Value = 25
$ws.Cells.Item(204, 18).Value = 'Hortaliza'
